$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Key input change: lower the discount rate from 6% to 4% ---
$ws2.Range("Z19").Value = 0.04

# --- Re-enter the projection formulas across their rows so Excel groups
#     them into shared formulas (matches the saved shape of the sheet) ---
$ws2.Range("N10:W10").Formula = '=M32*$Z$19'
$ws2.Range("B21:G21").Formula = '=B4/B2-1'

# --- Number-format tweak: show the price inputs/outputs with 2 decimals ---
$ws1.Range("D2").NumberFormat = "#,##0.00"
$ws2.Range("Z23").NumberFormat = "#,##0.00"

# --- Restore view/selection state ---
$ws1.Range("B15").Select()

$ws2.Activate()
$ws2.Range("Z23").Select()
